$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" row entirely (row 26). This shifts all following rows up by one.
$ws.Rows("26:26").Delete()

# After the shift above, the row that held "SC 92" is now row 27. Remove it too,
# shifting the remaining rows up by one more.
$ws.Rows("27:27").Delete()

# At this point rows 26-33 hold (in order): SC 5, SC 101, SC 105, SC 119, SC 120,
# SC 132, SC 193, SC 232 - matching the target layout. Now update the F column
# (error) values that changed as part of the re-run of the missing-data process.

# SC 5 (row 26): F value is now missing.
$ws.Range("F26").Value = ""

# SC 101 (row 27): F value is now present.
$ws.Range("F27").Value = 17

# SC 119 (row 29): F value is now missing.
$ws.Range("F29").Value = ""
